# Auto-generated edit script: update cryptos price/volume table (and swap RenderToken/EnergySwap rows)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Pre-format price cells whose new values look like plain numbers as Text,
# so they keep their original text ("0.9998" etc.) instead of becoming numeric.
$textCells = @("D5", "D6", "D7", "D8", "D9", "D11", "D13", "D14", "D15", "D16", "D18", "D20", "D22", "D23", "D24", "D25", "D26", "D27", "D29", "D30", "D31", "D32", "D33", "D34", "D36", "D40", "D42", "D43", "D44", "D45", "D46", "D47", "D48", "D49", "D50", "D51")
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# Apply the updated values
$ws.Range("D2").Value = "29.351.18"
$ws.Range("E2").Value = "  -0.09%  "
$ws.Range("D3").Value = "1.842.57"
$ws.Range("E3").Value = "  -0.25%  "
$ws.Range("E4").Value = "  -0.03%  "
$ws.Range("D5").Value = "239.22"
$ws.Range("E5").Value = "  -0.50%  "
$ws.Range("D6").Value = "0.6287"
$ws.Range("E6").Value = "  -0.40%  "
$ws.Range("D7").Value = "0.9998"
$ws.Range("D8").Value = "0.07524"
$ws.Range("E8").Value = "  -0.45%  "
$ws.Range("D9").Value = "0.2930"
$ws.Range("E9").Value = "  -1.07%  "
$ws.Range("E10").Value = "  -0.43%  "
$ws.Range("D11").Value = "0.07695"
$ws.Range("E11").Value = "  -0.37%  "
$ws.Range("D12").Value = "1.854.95"
$ws.Range("E12").Value = "  -6.55%  "
$ws.Range("D13").Value = "4.985"
$ws.Range("E13").Value = "  +0.11%  "
$ws.Range("D14").Value = "0.6775"
$ws.Range("E14").Value = "  -0.91%  "
$ws.Range("D15").Value = "0.00001042"
$ws.Range("E15").Value = "  +4.62%  "
$ws.Range("D16").Value = "82.82"
$ws.Range("E16").Value = "  +0.17%  "
$ws.Range("D17").Value = "2.084.94"
$ws.Range("E17").Value = "  -7.92%  "
$ws.Range("D18").Value = "6.098"
$ws.Range("E18").Value = "  -1.27%  "
$ws.Range("D19").Value = "29.367.80"
$ws.Range("D20").Value = "227.34"
$ws.Range("E20").Value = "  -1.59%  "
$ws.Range("E21").Value = "  -0.69%  "
$ws.Range("D22").Value = "0.9997"
$ws.Range("E22").Value = "  -0.03%  "
$ws.Range("D23").Value = "7.412"
$ws.Range("E23").Value = "  -2.22%  "
$ws.Range("D24").Value = "1.001"
$ws.Range("E24").Value = "  +0.01%  "
$ws.Range("D25").Value = "156.61"
$ws.Range("E25").Value = "  +1.27%  "
$ws.Range("D26").Value = "0.1384"
$ws.Range("E26").Value = "  -0.35%  "
$ws.Range("D27").Value = "8.347"
$ws.Range("E27").Value = "  -0.95%  "
$ws.Range("E28").Value = "  -0.42%  "
$ws.Range("D29").Value = "1.458"
$ws.Range("E29").Value = "  -0.98%  "
$ws.Range("D30").Value = "1.276"
$ws.Range("E30").Value = "  +1.22%  "
$ws.Range("D31").Value = "0.05632"
$ws.Range("E31").Value = "  -2.91%  "
$ws.Range("D32").Value = "4.093"
$ws.Range("E32").Value = "  -0.55%  "
$ws.Range("D33").Value = "4.015"
$ws.Range("E33").Value = "  +0.01%  "
$ws.Range("D34").Value = "1.833"
$ws.Range("E34").Value = "  -1.51%  "
$ws.Range("E35").Value = "  -0.29%  "
$ws.Range("D36").Value = "0.7063"
$ws.Range("E36").Value = "  -1.43%  "
$ws.Range("E37").Value = "  -0.29%  "
$ws.Range("D38").Value = "1.239.17"
$ws.Range("E39").Value = "  +0.02%  "
$ws.Range("D40").Value = "2.758"
$ws.Range("E40").Value = "  -1.17%  "
$ws.Range("D42").Value = "0.8987"
$ws.Range("E42").Value = "  -0.75%  "
$ws.Range("D43").Value = "0.9991"
$ws.Range("E43").Value = "  -0.08%  "
$ws.Range("D44").Value = "101.81"
$ws.Range("E44").Value = "  +0.49%  "
$ws.Range("D45").Value = "65.37"
$ws.Range("E45").Value = "  -2.55%  "
$ws.Range("D46").Value = "0.00000000119"
$ws.Range("E46").Value = "  +0.54%  "
$ws.Range("D47").Value = "7.016"
$ws.Range("E47").Value = "  -3.91%  "
$ws.Range("D48").Value = "0.3990"
$ws.Range("E48").Value = "  -0.39%  "
$ws.Range("B49").Value = "RenderToken"
$ws.Range("C49").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D49").Value = "1.666"
$ws.Range("E49").Value = "  -1.40%  "
$ws.Range("B50").Value = "EnergySwap"
$ws.Range("C50").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D50").Value = "8.885"
$ws.Range("E50").Value = "  -2.75%  "
$ws.Range("D51").Value = "0.1120"
$ws.Range("E51").Value = "  -0.09%  "
